$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# The category headers were reshuffled (new categories inserted
# alphabetically) and three brand-new trailing columns (I, J, K) were added.
# First, give the three new columns the same "header" style (bold + border)
# as the existing header cells, then set the text for every header cell.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "C INDUSTRIE"
$ws.Range("C1").Value = "E DISTRIBUTIE VAN WATER; AFVAL-EN AFVALWATERBEHEER EN SANERING "
$ws.Range("D1").Value = "G GROOT-EN DETAILHANDEL; REPARATIE VAN AUTO'S EN MOTORFIETSEN "
$ws.Range("E1").Value = "K FINANCIËLE ACTIVITEITEN EN VERZEKERINGEN "
$ws.Range("F1").Value = "M VRIJE BEROEPEN EN WETENSCHAPPELIJKE EN TECHNISCHE ACTIVITEITEN "
$ws.Range("G1").Value = "N ADMINISTRATIEVE EN ONDERSTEUNENDE DIENSTEN "
$ws.Range("H1").Value = "O OPENBAAR BESTUUR EN DEFENSIE; VERPLICHTE SOCIALE VERZEKERINGEN "
$ws.Range("I1").Value = "Other"
$ws.Range("J1").Value = "ROUTE INZAMELING"
$ws.Range("K1").Value = "S OVERIGE DIENSTEN "

# --- New data for columns I, J, K (rows 2-7) ---
$newData = @{
    2 = @{ I = 120049.539;            J = 529654.2749999999;  K = 43466.94300000002 }
    3 = @{ I = 99173.14699999998;     J = 538053.1150000003;  K = 43079.79 }
    4 = @{ I = 123341.016;            J = 520607.6140000001;  K = 41367.22300000001 }
    5 = @{ I = 136765.695;            J = 526426.5839999998;  K = 41816.752 }
    6 = @{ I = 129241.683;            J = 485015.3290000001;  K = 39440.93300000001 }
    7 = @{ I = 139829.306;            J = 392233.8750000002;  K = 38063.21500000001 }
}

foreach ($r in $newData.Keys) {
    $row = $newData[$r]
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
}

# --- Updated values for existing columns B..H (rows 2-7) ---
$updatedData = @{
    2 = @{ B = 34745.231;             C = 748006.3029999996;  D = 104158.087;        E = 57339.679;           F = 40664.467;           G = 34740.71900000002;  H = 128567.495 }
    3 = @{ B = 25664.922;             C = 737170.653;         D = 82913.72200000002; E = 59962.99900000001;   F = 50404.33300000002;   G = 32511.28399999999;  H = 157835.9609999999 }
    4 = @{ B = 28879.662;             C = 720690.9629999996;  D = 79585.27499999992; E = 51864.666;           F = 52235.398;           G = 32483.50600000002;  H = 174839.7639999999 }
    5 = @{ B = 35099.97099999998;     C = 700839.1999999998;  D = 81785.68800000002; E = 46912.32600000002;   F = 56226.63100000002;   G = 34692.66899999999;  H = 168733.342 }
    6 = @{ B = 48447.622;             C = 695208.6309999995;  D = 82534.40599999999; E = 56057.27100000003;   F = 71084.997;           G = 36254.30000000001;  H = 195480.117 }
    7 = @{ B = 41879.069;             C = 804773.96;          D = 76443.864;         E = 54761.91900000003;   F = 79418.52300000004;   G = 39276.80599999998;  H = 291652.6330000002 }
}

foreach ($r in $updatedData.Keys) {
    $row = $updatedData[$r]
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
}

